$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$aw = $excel.ActiveWindow

function TryProp($name, $val) {
    try {
        $aw.$name = $val
        Write-Host "SET $name = $val -> OK, read back:" $aw.$name
    } catch {
        Write-Host "SET $name FAILED:" $_
    }
}

TryProp "SplitRow" 1
TryProp "SplitColumn" 0
TryProp "ScrollRow" 6
TryProp "ScrollColumn" 1
